$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-unused Col2/Col3 columns (D and E) entirely.
$ws.Range("D1:E3").Clear()

# Row 3 ("TC02") loses its TempTC02 / tc02col1 data; the name in B3 goes away,
# and C3 is reused to hold the CMS1500_1.png filename (given a distinct font).
$ws.Range("B3").Clear()
$ws.Range("C3").Value = "CMS1500_1.png"
$ws.Range("C3").Font.Bold = $true

# Match the final selection recorded in the sheet view.
$ws.Range("D3").Select() | Out-Null
